$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new column before the old "F" (type) column.
#    Everything from F onward (type, amount, date, userId) shifts
#    right by one: F->G, G->H, H->I, I->J.
# ------------------------------------------------------------------
$ws.Columns("F:F").Insert()

# Header for the newly inserted column
$ws.Range("F1").Value = "Unnamed: 0.1.1.1.1"

# ------------------------------------------------------------------
# 2) Fill the new column F (rows 2-17) with the same running index
#    that already lives in column E for those rows (0..15).
# ------------------------------------------------------------------
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 6).Value = $r - 2
}

# ------------------------------------------------------------------
# 3) Fix the "staircase" of index columns for the tail rows, now
#    that there are six index columns (A-F) instead of five (A-E).
#    Row 18: A-E filled with 16, F stays blank.
#    Row 19: A-C filled with 17, D-F stay blank.
#    Row 20: A-B filled with 18, C-F stay blank.
# ------------------------------------------------------------------
$ws.Range("E18").Value = 16
$ws.Range("C19").Value = 17
$ws.Range("B20").Value = 18

# ------------------------------------------------------------------
# 4) Data fix: row 5's amount (entertain) changes from 80 to 70.
#    After the column insert, amount now lives in column H.
# ------------------------------------------------------------------
$ws.Range("H5").Value = 70

# ------------------------------------------------------------------
# 5) Append a brand-new record as row 21 (index 19).
# ------------------------------------------------------------------
$ws.Range("A20").Copy() | Out-Null
$ws.Range("A21").PasteSpecial(-4122) | Out-Null

$ws.Range("A21").Value = 19
$ws.Range("G21").Value = "clothing"
$ws.Range("H21").Value = 40

# Force the date-looking text to stay a literal string (matches the
# other "date" column cells, which are all stored as text, not real
# Excel dates).
$ws.Range("I21").NumberFormat = "@"
$ws.Range("I21").Value = "2021-08-20"
$ws.Range("I21").Style = "Normal"

$ws.Range("J21").Value = "Dixon3220"
